# Add a new table-of-contents entry for "Kosuri et al" above the current
# row 21 ("Gao et al"), shifting all subsequent entries down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 21 (pushes existing rows 21-47 down to 22-48)
$ws.Rows("21:21").Insert()

# Populate the new row. Order matches how the entry was originally authored
# (title first, then author, doi, association, comments).
$ws.Range("B21").Value = "Composability of regulatory sequences controlling transcription and translation in Escherichia coli"
$ws.Range("A21").Value = "Kosuri et al"
$ws.Range("C21").Value = "10.1073/pnas.1301301110"
$ws.Range("D21").Value = "Promoter engineering"
$ws.Range("E21").Value = "Kosuri et al. synthesized 12,563 combinations of common promoters and ribosome binding sites and simultaneously measured DNA, RNA, and protein levels from the entire library. "

# Reflect the cursor/selection left at E21 after data entry.
$ws.Range("E21").Select() | Out-Null
